$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1804.1111
$ws.Range("I28").Value = 1954.625
$ws.Range("K28").Value = 1954.625
$ws.Range("M28").Value = -1469.625
$ws.Range("H86").Value = 5788.4614
$ws.Range("I86").Value = 3929
$ws.Range("K86").Value = 3929
$ws.Range("M86").Value = -2806
$ws.Range("H89").Value = 5788.4614
$ws.Range("I89").Value = 3929
$ws.Range("K89").Value = 19645
$ws.Range("M89").Value = -14029
$ws.Range("H98").Value = 3813.2144
$ws.Range("I98").Value = 1048
$ws.Range("K98").Value = 1048
$ws.Range("M98").Value = 450
$ws.Range("H118").Value = 1351
$ws.Range("I118").Value = 1351
$ws.Range("K118").Value = 4053
$ws.Range("M118").Value = -2396
$ws.Range("H122").Value = 3813.2144
$ws.Range("I122").Value = 1048
$ws.Range("K122").Value = 3144
$ws.Range("M122").Value = -694
$ws.Range("H137").Value = 4000
$ws.Range("I137").Value = 4000
$ws.Range("K137").Value = 12000
$ws.Range("M137").Value = -9450

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1999.8889
$ws.Range("I45").Value = 1999.875
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1999.875
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1622.875
$ws.Range("N45").Value = -2754
$ws.Range("H69").Value = 332500
$ws.Range("J69").Value = 332500
$ws.Range("L69").Value = 332500
$ws.Range("N69").Value = -333998
$ws.Range("H72").Value = 332500
$ws.Range("J72").Value = 332500
$ws.Range("L72").Value = 997500
$ws.Range("N72").Value = -1004988
$ws.Range("H92").Value = 65598.5
$ws.Range("J92").Value = 65598.5
$ws.Range("L92").Value = 65598.5
$ws.Range("N92").Value = -70590.5
$ws.Range("H122").Value = 1416.64
$ws.Range("I122").Value = 1013.35
$ws.Range("J122").Value = 3029.8
$ws.Range("K122").Value = 3040.05
$ws.Range("L122").Value = 9089.400000000001
$ws.Range("M122").Value = -590.0500000000002
$ws.Range("N122").Value = -13989.4
$ws.Range("H127").Value = 40000
$ws.Range("I127").Value = 40000
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 40000
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -35040
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 1917.0834
$ws.Range("I132").Value = 1782.4348
$ws.Range("K132").Value = 5347.3044
$ws.Range("M132").Value = -2817.3044

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 31990
$ws.Range("J6").Value = 31990
$ws.Range("L6").Value = 31990
$ws.Range("N6").Value = -32216
$ws.Range("H20").Value = 2754
$ws.Range("I20").Value = 2500
$ws.Range("K20").Value = 2500
$ws.Range("M20").Value = -2253
$ws.Range("H22").Value = 532.8
$ws.Range("I22").Value = 616
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 616
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -443
$ws.Range("N22").Value = -546
$ws.Range("H37").Value = 4544.8
$ws.Range("J37").Value = 4544.8
$ws.Range("L37").Value = 4544.8
$ws.Range("N37").Value = -4818.8
$ws.Range("H92").Value = 65000
$ws.Range("J92").Value = 65000
$ws.Range("L92").Value = 65000
$ws.Range("N92").Value = -69992
$ws.Range("H102").Value = 2956
$ws.Range("I102").Value = 2956
$ws.Range("K102").Value = 2956
$ws.Range("M102").Value = 289
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 383.16666
$ws.Range("I10").Value = 199.75
$ws.Range("K10").Value = 199.75
$ws.Range("M10").Value = -60.75
$ws.Range("H22").Value = 85001.336
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 85001.336
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 85001.336
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -85701.336
$ws.Range("H132").Value = 1507.9375
$ws.Range("I132").Value = 1215.8334
$ws.Range("K132").Value = 3647.5002
$ws.Range("M132").Value = -1117.5002
$ws.Range("H134").Value = 1941.6774
$ws.Range("I134").Value = 1875.7778
$ws.Range("K134").Value = 5627.3334
$ws.Range("M134").Value = -3092.3334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 19314
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 19314
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 57942
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -58482
$ws.Range("H67").Value = 19314
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 19314
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 57942
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -59814
$ws.Range("H131").Value = 2412.2104
$ws.Range("I131").Value = 1788.8
$ws.Range("K131").Value = 5366.4
$ws.Range("M131").Value = -326.3999999999996

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4873.5
$ws.Range("I70").Value = 4799
$ws.Range("J70").Value = 4948
$ws.Range("K70").Value = 4799
$ws.Range("L70").Value = 4948
$ws.Range("M70").Value = -4529
$ws.Range("N70").Value = -5488
$ws.Range("H73").Value = 4873.5
$ws.Range("I73").Value = 4799
$ws.Range("J73").Value = 4948
$ws.Range("K73").Value = 4799
$ws.Range("L73").Value = 4948
$ws.Range("M73").Value = -3863
$ws.Range("N73").Value = -6820
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H113").Value = 1532.2
$ws.Range("I113").Value = 1602.4445
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 1602.4445
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = 567.5554999999999
$ws.Range("N113").Value = -5240

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2751
$ws.Range("I46").Value = 1035.8889
$ws.Range("J46").Value = 3780.0667
$ws.Range("K46").Value = 1035.8889
$ws.Range("L46").Value = 3780.0667
$ws.Range("M46").Value = -847.8888999999999
$ws.Range("N46").Value = -4156.066699999999
$ws.Range("H61").Value = 3922.4666
$ws.Range("I61").Value = 3746.0454
$ws.Range("J61").Value = 4407.625
$ws.Range("K61").Value = 3746.0454
$ws.Range("L61").Value = 4407.625
$ws.Range("M61").Value = -3544.0454
$ws.Range("N61").Value = -4811.625
$ws.Range("H101").Value = 19000
$ws.Range("J101").Value = 19000
$ws.Range("L101").Value = 19000
$ws.Range("N101").Value = -25490
$ws.Range("H113").Value = 3922.4666
$ws.Range("I113").Value = 3746.0454
$ws.Range("J113").Value = 4407.625
$ws.Range("K113").Value = 3746.0454
$ws.Range("L113").Value = 4407.625
$ws.Range("M113").Value = -1576.0454
$ws.Range("N113").Value = -8747.625
$ws.Range("H122").Value = 6762.5625
$ws.Range("I122").Value = 6584.2173
$ws.Range("K122").Value = 19752.6519
$ws.Range("M122").Value = -17302.6519

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 276.7143
$ws.Range("I113").Value = 137.4
$ws.Range("K113").Value = 412.2
$ws.Range("M113").Value = 1757.8
$ws.Range("H136").Value = 4341.85
$ws.Range("I136").Value = 4367.0586
$ws.Range("J136").Value = 4199
$ws.Range("K136").Value = 13101.1758
$ws.Range("L136").Value = 12597
$ws.Range("M136").Value = -10551.1758
$ws.Range("N136").Value = -17697
